# RDCC-5182 Added Version check
#
# Inserts a new "VERSION" worksheet between "Staff Data" and "Sheet2" that
# contains a small "File version" / "vx.xx" key-value pair in row 6, and
# makes it the active/selected tab (mirroring the author's manual edit in
# Excel: a new sheet was added, the old "tabSelected" flag moved off
# "Staff Data" onto the new sheet, and that new sheet became the active
# tab of the workbook).

$wb = $excel.ActiveWorkbook

$staffData = $wb.Sheets.Item("Staff Data")
$sheet2    = $wb.Sheets.Item("Sheet2")

# Copy the existing "Sheet2" so the new sheet inherits the same sheet
# formatting (row height, etc.) as its siblings, then drop into place
# right after "Staff Data" (i.e. before "Sheet2") — matching the sheet
# order in the workbook: Staff Data, VERSION, Sheet2.
$sheet2.Copy($null, $staffData) | Out-Null
$newSheet = $wb.Sheets.Item("Sheet2 (2)")

# Start from a clean sheet (remove the copied header row/content).
$newSheet.Cells.Clear() | Out-Null
$newSheet.Rows.Item(1).Delete() | Out-Null

$newSheet.Name = "VERSION"

$newSheet.Range("A6").Value = "File version"
$newSheet.Range("B6").Value = "vx.xx"

# Make VERSION the active/selected sheet, with B6 as the selected cell —
# this is what toggles tabSelected onto this sheet (and off "Staff Data")
# and sets the workbook's activeTab.
$newSheet.Activate() | Out-Null
$newSheet.Range("B6").Select() | Out-Null
